$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks fully numeric need an explicit
# text format, otherwise Excel auto-converts the assigned string into
# a Number cell (losing the intended plain-text "Price" representation).
$textProtectRefs = @(
    'D5', 'D6', 'D9', 'D11', 'D12', 'D14', 'D19', 'D20',
    'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29',
    'D31', 'D33', 'D34', 'D35', 'D36', 'D38', 'D39', 'D40',
    'D41', 'D43', 'D45', 'D47', 'D48', 'D51'
)
foreach ($ref in $textProtectRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.402.91'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '3.139.32'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '609.00'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").Value = '143.64'
$ws.Range("E6").Value = '  -1.97%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.139.83'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("E13").Value = '  +3.27%  '
$ws.Range("D14").Value = '35.37'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '3.657.59'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").Value = '64.371.97'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '3.140.72'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("D19").Value = '6.85'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '477.01'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '14.83'
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("D22").Value = '0.718'
$ws.Range("E22").Value = '  +2.03%  '
$ws.Range("D23").Value = '7.74'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '85.52'
$ws.Range("E24").Value = '  +3.40%  '
$ws.Range("D25").Value = '13.40'
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  -3.41%  '
$ws.Range("D28").Value = '8.42'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = '7.24'
$ws.Range("E29").Value = '  +7.60%  '
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("D31").Value = '2.05'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").Value = '26.81'
$ws.Range("E33").Value = '  +3.10%  '
$ws.Range("D34").Value = '2.63'
$ws.Range("E34").Value = '  -3.71%  '
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").Value = '  +0.62%  '
$ws.Range("D36").Value = '5.97'
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").Value = '0.0₃0759'
$ws.Range("E37").Value = '  +4.31%  '
$ws.Range("D38").Value = '52.54'
$ws.Range("E38").Value = '  -2.00%  '
$ws.Range("D39").Value = '3.03'
$ws.Range("E39").Value = '  +3.72%  '
$ws.Range("D40").Value = '445.87'
$ws.Range("E40").Value = '  -3.28%  '
$ws.Range("D41").Value = '0.0394'
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("D43").Value = '8.25'
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("D44").Value = '2.888.54'
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("D45").Value = '0.261'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").Value = '2.39'
$ws.Range("E47").Value = '  +3.69%  '
$ws.Range("D48").Value = '26.21'
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  -0.51%  '
$ws.Range("D51").Value = '119.63'
$ws.Range("E51").Value = '  +0.56%  '

